# bga-obra.xlsx restructuring:
#   Hoja1 (mixed obra+referente data) -> renamed "bga-obra" (data kept as-is)
#   + new "obra" sheet (ID, Titulo, Fecha, Dimensiones, Tecnica, Archivo)
#   + new "referente" sheet (ID, Titulo, Fecha, Periodico, Archivo)

$wb = $excel.ActiveWorkbook
$bga = $wb.Worksheets.Item(1)

# --- Build "obra" sheet as a copy of the original sheet, so it inherits the
# --- exact same column widths/formatting, then strip it down to just the
# --- "obra" shaped data. Column F (Periodico) is empty/has no custom width
# --- in the source sheet, so deleting it cleanly shifts G (Archivo, 25.5)
# --- into F with no stray width entries left behind.
$bga.Copy([Type]::Missing, $bga)
$obra = $wb.Worksheets.Item(2)
$obra.Name = "obra"
$obra.Columns.Item(6).Delete()
$obra.Cells.ClearContents()

$obra.Range("A1").Value = "ID"
$obra.Range("B1").Value = "Título"
$obra.Range("C1").Value = "Fecha"
$obra.Range("D1").Value = "Dimensiones"
$obra.Range("E1").Value = "Técnica"
$obra.Range("F1").Value = "Archivo"

$obra.Range("A2").Value = 1
$obra.Range("B2").Value = "Los Suicidas del Sisga No 1"
$obra.Range("C2").Value = 1965
$obra.Range("D2").Value = "120x100cm"
$obra.Range("E2").Value = "Óleo sobre lienzo"
$obra.Range("F2").Value = "suicidas-sisga-1.jpg"

$obra.Range("A3").Value = 2
$obra.Range("B3").Value = "El Paraíso"
$obra.Range("C3").Value = 1997
$obra.Range("D3").Value = "160x45 cm"
$obra.Range("E3").Value = "Óleo sobre lienzo"
$obra.Range("F3").Value = "el-paraiso.jpg"

$obra.Range("A4").Value = 3
$obra.Range("B4").Value = "Zócalo de la tragedia"
$obra.Range("C4").Value = 1983
$obra.Range("D4").Value = "100x71"
$obra.Range("E4").Value = "Tipografía sobre papel"
$obra.Range("F4").Value = "zocalo-tragedia.jpg"

# --- Build "referente" sheet the same way: copy the original sheet (for
# --- matching formatting), drop the obra-only columns (Dimensiones,
# --- Tecnica), then repopulate with the referente-shaped rows.
$bga.Copy([Type]::Missing, $obra)
$referente = $wb.Worksheets.Item(3)
$referente.Name = "referente"
$referente.Columns("D:E").Delete()
$referente.Cells.ClearContents()

$referente.Range("A1").Value = "ID"
$referente.Range("B1").Value = "Título"
$referente.Range("C1").Value = "Fecha"
$referente.Range("D1").Value = "Periódico"
$referente.Range("E1").Value = "Archivo"

$referente.Range("A2").Value = 1
$referente.Range("B2").Value = "Doble suicidio en ""El Sisga"""
$referente.Range("C2").Value = "Junio 29 1965"
$referente.Range("D2").Value = "El Tiempo"
$referente.Range("E2").Value = "doble-suicidio-el-tiempo.jpg"

$referente.Range("A3").Value = 2
$referente.Range("B3").Value = "Una indígena y su hijo murieron en persecución"
$referente.Range("C3").Value = "Mayo 24 del 96"
$referente.Range("D3").Value = "El Tiempo"
$referente.Range("E3").Value = "indigena-hijo-el-tiempo.jpg"

$referente.Range("A4").Value = 3
$referente.Range("B4").Value = "Láminas de paisajes latinoamericanos"
$referente.Range("E4").Value = "laminas-paisajes.jpg"

$referente.Range("A5").Value = 4
$referente.Range("B5").Value = "Exmilitar Mata a la Esposa de su Amigo y se Suicida"
$referente.Range("E5").Value = "exmilitar-mata-esposa.jpg"

# Column widths that land on the 1/6-character grid this engine quantizes
# to can be set exactly; reuse them where the copied width isn't already
# correct.
$referente.Columns.Item(5).ColumnWidth = 24.666666666666668
$referente.Columns.Item(2).ColumnWidth = 34.998697916666664

# --- Rename the original sheet and clear its "current" selection/active
# --- styling now that it is no longer the freshly-opened tab.
$bga.Name = "bga-obra"

# --- View state: "obra" becomes the active/selected tab; "bga-obra" selects
# --- the whole header row; "referente" keeps an arbitrary out-of-range cell
# --- selected, matching the source workbook's saved state.
$bga.Activate()
$bga.Rows.Item(1).Select()

$referente.Activate()
$referente.Range("E8").Select()

$obra.Activate()
$obra.Range("A5").Select()
